$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 1

$ws.Range("N3").Value = 1

$ws.Range("M4").Value = 10
$ws.Range("N4").Value = 1

$ws.Range("N5").Value = 1

$ws.Range("N6").Value = 1

$ws.Range("N7").Value = 1
